$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (AT trial header values) - update B1:E1
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON, meanEMG legmaxROM) - update B2:E2
$ws.Range("B2").Value = 67.344776313292996
$ws.Range("C2").Value = 37.31198412568245
$ws.Range("D2").Value = 67.458618079301345
$ws.Range("E2").Value = 40.905546271713845

# Row 3 (STR, meanEMG legmaxROM) - update B3:E3
$ws.Range("B3").Value = 63.376823674849284
$ws.Range("C3").Value = 43.216688876332171
$ws.Range("D3").Value = 43.216688876332171
$ws.Range("E3").Value = 54.251121620335695

# Update the selection to match the new active selection range B1:E3
$ws.Range("B1:E3").Select()
